# The default wait time used to be specified in milliseconds (e.g. 7000 / 6000)
# but WebDriverWait actually expects seconds, so all timeout-style values are
# being corrected to their second-based equivalents, and a clarifying comment
# is added wherever a timeout value is configured.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c-demo_ui")
$ws.Activate()

# Fix the "Plateform"/"plateform" typo -> "Platform"/"platform" in the
# comments column for every platform/version pair.
$ws.Range("C7").Value  = "Platform type for remote web driver intializing"
$ws.Range("C8").Value  = "Version for platform type selected"
$ws.Range("C9").Value  = "Platform type for remote web driver intializing"
$ws.Range("C10").Value = "Version for platform type selected"
$ws.Range("C11").Value = "Platform type for remote web driver intializing"
$ws.Range("C12").Value = "Version for platform type selected"
$ws.Range("C13").Value = "Platform type for remote web driver intializing"
$ws.Range("C14").Value = "Version for platform type selected"

# waitTime was 7000 (ms); change to 7 (seconds) and document why.
$ws.Range("B23").Value = "7"
$ws.Range("C23").Value = "Wait time delay is seconds, not milliseconds"

$ws.Range("C16").Select() | Out-Null
